# "Coreação da Agenda" - update status column (F) on the project agenda sheet,
# and move the active selection/view as the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row that was "?" (Em aberto) and is now "*" (Duvida) - write this first so
# the new shared string "*" lands before "!!" in the shared-strings table.
$ws.Range("F32").Value = "*"

# Rows that were "!" (Concluida) and are now "!!" (Revisado ok)
$ws.Range("F4").Value = "!!"
$ws.Range("F5").Value = "!!"
$ws.Range("F6").Value = "!!"

# Rows that were blank and now get "!!" (Revisado ok)
$ws.Range("F7").Value = "!!"
$ws.Range("F9").Value = "!!"
$ws.Range("F11").Value = "!!"
$ws.Range("F13").Value = "!!"
$ws.Range("F15").Value = "!!"

# Rows that were "?" (Em aberto) and are now "!" (Concluida)
$ws.Range("F27").Value = "!"
$ws.Range("F28").Value = "!"
$ws.Range("F30").Value = "!"
$ws.Range("F33").Value = "!"
$ws.Range("F36").Value = "!"
$ws.Range("F38").Value = "!"
$ws.Range("F40").Value = "!"

# Leave the view scrolled and selection where the author left off
$ws.Range("H8").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
